# This script applies the update described by the commit:
# - A new claim (Caso 5979) is inserted as row 16, pushing the existing
#   rows 16-50 down by one (now rows 17-51).
# - Five brand-new claims are appended at the bottom of the table
#   (rows 52-55, plus the former last row which is now row 51).
#
# All "text-looking" columns (A-L) in this sheet are stored as text, even
# when they contain only digits (case numbers, OT numbers, comuna, etc.),
# so we force that with a leading single-quote the same way a user typing
# into Excel would. Columns M/N (coordinates) are genuine numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $text)
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

function Set-NumberCell {
    param($row, $col, $number)
    $ws.Cells.Item($row, $col).Value = $number
}

function Set-BlankCell {
    param($row, $col)
    # Touch a formatting property (no visual change) so the engine keeps
    # an explicit, empty cell entry instead of dropping it altogether.
    $ws.Cells.Item($row, $col).Font.Bold = $false
}

# ---------------------------------------------------------------------
# 1) Insert the new row 16 (shifts old rows 16..50 down to 17..51)
# ---------------------------------------------------------------------
$ws.Rows.Item(16).Insert()

Set-TextCell 16 1 "5979"
Set-TextCell 16 2 "2/24/2025"
Set-TextCell 16 3 "CHILAVERT, MARTINIANO, CORONEL 3114"
Set-TextCell 16 4 "8"
Set-TextCell 16 5 "803608474"
Set-TextCell 16 6 "Optical Power"
Set-TextCell 16 7 "Pendiente"
Set-TextCell 16 8 "Picada"
Set-TextCell 16 9 "1"
Set-TextCell 16 10 "Cambio"
Set-TextCell 16 11 "Sin equipos"
Set-TextCell 16 12 "Pasante"
Set-NumberCell 16 13 -58.444984
Set-NumberCell 16 14 -34.659504

# ---------------------------------------------------------------------
# 2) Append five new claims at the end of the table (rows 52-55)
# ---------------------------------------------------------------------

# Row 52 - Caso 3715
Set-TextCell 52 1 "3715"
Set-TextCell 52 2 "6/4/2025"
Set-TextCell 52 3 "EL SERENO 358"
Set-TextCell 52 4 "10"
Set-TextCell 52 5 "807168098"
Set-TextCell 52 6 "Optical Power"
Set-TextCell 52 7 "Pendiente"
Set-TextCell 52 8 "Poste inclinado"
Set-TextCell 52 9 "1"
Set-TextCell 52 10 "Aplomo"
Set-TextCell 52 11 "Sin equipos"
Set-TextCell 52 12 "Poste"
Set-NumberCell 52 13 -58.487371
Set-NumberCell 52 14 -34.640099

# Row 53 - Caso 5997
Set-TextCell 53 1 "5997"
Set-TextCell 53 2 "6/4/2025"
Set-TextCell 53 3 "MARMOL, JOSE 256"
Set-TextCell 53 4 "5"
Set-TextCell 53 5 "807187768"
Set-TextCell 53 6 "Optical Power"
Set-TextCell 53 7 "Pendiente"
Set-TextCell 53 8 "Picada coincide con reclamo de cables con mismo numero de caso"
Set-TextCell 53 9 "1"
Set-TextCell 53 10 "Cambio"
Set-TextCell 53 11 "Sin equipos"
Set-TextCell 53 12 "Pasante"
Set-NumberCell 53 13 -58.425845
Set-NumberCell 53 14 -34.616562

# Row 54 - Caso 807187860 (note trailing space kept on purpose in column A)
Set-TextCell 54 1 "807187860 "
Set-TextCell 54 2 "6/4/2025"
Set-TextCell 54 3 "Av. San Juan 3960"
Set-TextCell 54 4 "5"
Set-TextCell 54 5 "807187860"
Set-TextCell 54 6 "Optical Power"
Set-TextCell 54 7 "Pendiente"
Set-TextCell 54 8 "Colocar columna contactar a Matias Tapia 1171744701 por si hay alguna duda o problema que surja en el momento ya que es para posterior tendido de ftth"
Set-TextCell 54 9 "1"
Set-TextCell 54 10 "Cambio"
Set-TextCell 54 11 "Sin equipos"
Set-TextCell 54 12 "Pasante"
Set-BlankCell 54 13
Set-BlankCell 54 14

# Row 55 - Caso 807187874
Set-TextCell 55 1 "807187874"
Set-TextCell 55 2 "6/4/2025"
Set-TextCell 55 3 "Corvalan 4348"
Set-TextCell 55 4 "8"
Set-TextCell 55 5 "807187874"
Set-TextCell 55 6 "Optical Power"
Set-TextCell 55 7 "Pendiente"
Set-TextCell 55 8 "Colocar columna barrio papa francis contactar a inspector Matias Tapia 1171744701 para el ingreso al barrio ya que el esta en obra en el lugar"
Set-TextCell 55 9 "1"
Set-TextCell 55 10 "Cambio"
Set-TextCell 55 11 "Sin equipos"
Set-TextCell 55 12 "Pasante"
Set-NumberCell 55 13 -58.462456
Set-NumberCell 55 14 -34.674961

Write-Output "done"
